# fix: fix overflow title
# The "Chapter horizontal" and "Chapter vertical" slide layouts each ship
# with a placeholder ("Text Placeholder 3", the ctrTitle chapter-number
# box) whose sample text is the two-digit string "00". That default
# overflowed/looked wrong, so bump it to "01".
$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$layouts = $master.CustomLayouts

for ($i = 1; $i -le $layouts.Count; $i++) {
    $layout = $layouts.Item($i)
    $shapes = $layout.Shapes
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shape = $shapes.Item($j)
        if ($shape.HasTextFrame -eq -1) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "00") {
                $tr.Text = "01"
            }
        }
    }
}
